$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.150.73'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.677.33'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'214.22"
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('E6').Value = '  -4.42%  '
$ws.Range('D7').Value = "'1.007"
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.2677"
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = "'0.06340"
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').Value = "'21.24"
$ws.Range('E10').Value = '  -3.62%  '
$ws.Range('D11').Value = "'0.07602"
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '1.698.05'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = "'4.507"
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').Value = "'0.5707"
$ws.Range('E14').Value = '  -1.62%  '
$ws.Range('D15').Value = "'0.000008184"
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').Value = "'65.98"
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').Value = '26.199.31'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = "'4.844"
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = "'10.68"
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').Value = "'188.80"
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').Value = "'6.204"
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'148.92"
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').Value = "'0.1254"
$ws.Range('E25').Value = '  -4.34%  '
$ws.Range('D26').Value = "'7.680"
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').Value = "'15.81"
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = "'0.06395"
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = "'1.362"
$ws.Range('E29').Value = '  -2.45%  '
$ws.Range('D30').Value = "'1.307"
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('D31').Value = "'3.550"
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').Value = "'3.536"
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = "'1.666"
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').Value = "'1.012"
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('D35').Value = "'0.6049"
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').Value = "'2.739"
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').Value = "'0.01631"
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('D39').Value = "'6.143"
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('D40').Value = '1.088.28'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').Value = "'0.8752"
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').Value = "'100.14"
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('D44').Value = '1.832.29'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = "'0.00000000109"
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'57.29"
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = "'0.9976"
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = "'8.047"
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('D49').Value = "'0.05258"
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').Value = "'0.4279"
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = "'5.971"
$ws.Range('E51').Value = '  -1.28%  '
